$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting rows 94:202 down to 95:203
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with its data
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44601
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = 100112052
$ws.Cells.Item(94, 7).Value = "Albahaca"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 45
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 6000
$ws.Cells.Item(94, 13).Value = 6000
$ws.Cells.Item(94, 14).Value = "`$/paquete"
$ws.Cells.Item(94, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(94, 16).Value = 6000
$ws.Cells.Item(94, 17).Value = 1
$ws.Cells.Item(94, 18).Value = "Hortaliza"
